$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.616.62"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  +3.82%  "
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.14"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.02%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.09"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("E7").Value = "  +0.01%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.08"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.369"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +3.21%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.00"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +8.64%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0762"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("E12").Value = "  +2.59%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.52"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +7.59%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.808"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +5.15%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.191.66"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +1.99%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.10"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +3.82%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.920.43"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +2.56%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.570.79"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +3.77%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.53"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +1.39%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0861"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +4.81%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "250.03"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +2.51%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.36"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +4.17%  "
$ws.Range("E23").Value = "  +3.58%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E26").Value = "  +0.45%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.03"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +1.68%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.82"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +2.91%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.71"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("E31").Value = "  +7.11%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0618"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +4.92%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.94"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +6.69%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.34"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +4.05%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0892"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +21.31%  "
$ws.Range("E36").Value = "  +0.15%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.52"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +3.54%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +3.51%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.76"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +51.09%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.04"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +5.00%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.88"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +11.07%  "
$ws.Range("E42").Value = "  +4.53%  "
$ws.Range("E43").Value = "  -0.45%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.87"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +20.03%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.11"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +3.21%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.342.93"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +2.94%  "
$ws.Range("E47").Value = "  -0.14%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0815"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("E50").Value = "  +2.08%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.35"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +2.64%  "
